# "third teams time variable" - reorder the team-name lists stored as
# Python-literal strings in column O ("tied_teams") so the ordering
# matches the other columns' convention.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$groups = @(
    @{ Rows = @(3, 4); Value = "['Hungary', 'Uruguay']" },
    @{ Rows = @(9, 10, 11, 12); Value = "['Northern Ireland', 'Bulgaria']" },
    @{ Rows = @(53, 54, 55, 56, 57, 58, 59, 60); Value = "['Ireland', 'Argentina', 'Colombia', 'Costa Rica']" },
    @{ Rows = @(61, 62); Value = "['Argentina', 'Colombia']" },
    @{ Rows = @(63, 64, 65, 66, 67, 68, 69, 70, 71, 72, 73); Value = "['Scotland', 'Austria', 'Argentina', 'Colombia']" },
    @{ Rows = @(82); Value = "['South Korea', 'Cameroon']" },
    @{ Rows = @(104); Value = "['Netherlands', 'United States']" }
)

foreach ($g in $groups) {
    foreach ($r in $g.Rows) {
        $ws.Range("O$r").Value = $g.Value
    }
}
